$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "f" values from E5:E6 (leftover std-filter marker),
# now that std filtering uses a separate method.
$ws.Range("E5:E6").ClearContents()

# Move the active selection to H7 (matches the author's saved cursor state).
$ws.Range("H7").Select()
